$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after "Monthly Trend" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match page margins used by the other sheets (0.75/0.75/1/1/0.5/0.5 in)
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Copy the bold/bordered/centered header style from an existing sheet,
# then overwrite the text, so the style index is reused (not duplicated).
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("A1:B1"))
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("C1:D1"))

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the date-formatted style (column A data cells) down the A column
# before filling in values, again so the style index is reused. The source
# is a single cell so Excel tiles its format across the whole destination.
$wsWeekly.Range("A2").Copy($wsForecast.Range("A2:A20"))

$rows = @(
    @(45039.99999999999, 30, 6.421621811809398, 53.87072141995785),
    @(45060.99999999999, 30, 7.013960185754557, 53.11130850165905),
    @(45067.99999999999, 30, 6.087656879646235, 51.88507453389086),
    @(45165.99999999999, 28, 3.009164626348861, 52.29737418674591),
    @(45305.99999999999, 25, 1.272437945678684, 47.900935356594),
    @(45354.99999999999, 24, 3.092910062252504, 47.01408316122905),
    @(45494.99999999999, 22, -4.316480137290442, 44.95097696667724),
    @(45522.99999999999, 21, -1.528574343632132, 47.54732285664771),
    @(45585.99999999999, 20, -2.358214107242608, 43.44485253872237),
    @(45592.99999999999, 20, -2.756269551712653, 42.59024756412403),
    @(45599.99999999999, 20, -2.778501426791852, 43.04795000308468),
    @(45606.99999999999, 20, -4.679678016378832, 42.57361481454888),
    @(45613.99999999999, 20, -3.154080038204632, 42.6354837039874),
    @(45620.99999999999, 20, -2.449765678550508, 43.10095254115106),
    @(45627.99999999999, 20, -1.990457479207063, 42.90020509669757),
    @(45634.99999999999, 19, -4.182519251029402, 43.69162757112688),
    @(45641.99999999999, 19, -2.593884983699295, 42.23153955995213),
    @(45648.99999999999, 19, -2.695125580233679, 42.3385643736105),
    @(45655.99999999999, 19, -4.210122873043985, 40.8141877927287)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Restore the original active sheet/tab selection
$wsWeekly.Activate()
[void]$wsWeekly.Range("A1").Select()
